$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H3: replace the placeholder comment with a real one
$ws.Range("H3").Value = "Déjà au local"

# Row 9 ("Clock 4 MHz - IQXO-70"): fill in pricing/quantity + dimensions,
# drop the leftover placeholder comment, and add the datasheet hyperlink.
$ws.Range("D9").Value = 1.96
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1.96
$ws.Range("H9").Value = ""

$ws.Hyperlinks.Add($ws.Range("N9"), "http://www.rs-particuliers.com/WebCatalog/Oscillateur_a_quartz__4_MHz__%C2%B125ppm_HCMOS__15pF__CMS__7_x_5_x_14mm__4_broches-6720814.aspx")
$ws.Range("N9").Style = "Hyperlink"

$ws.Range("K9").Value = "7*5*1,4"

# Move the active selection
$ws.Range("H16").Select()

Write-Host "done"
